# Weekly update: insert a new data row for "Terminal La Palmera de La Serena - Zanahoria"
# above current row 363, pushing the existing rows 363:457 down to 364:458.
# The new row duplicates the surrounding record's field values but carries the
# newest reporting date (serial 44932 = 2023-01-06).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 363 (shifts 363:457 -> 364:458,
# inherits number formatting from the row above, as Excel normally does).
$ws.Rows("363").Insert()

# Populate the newly inserted row 363.
$ws.Cells.Item(363, 1).Value = 8
$ws.Cells.Item(363, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(363, 3).Value = "Coquimbo"
$ws.Cells.Item(363, 4).Value = 44932
$ws.Cells.Item(363, 5).Value = 4
$ws.Cells.Item(363, 6).Value = 100114013
$ws.Cells.Item(363, 7).Value = "Zanahoria"
$ws.Cells.Item(363, 8).Value = "Sin especificar"
$ws.Cells.Item(363, 9).Value = "Primera"
$ws.Cells.Item(363, 10).Value = 600
$ws.Cells.Item(363, 11).Value = 5500
$ws.Cells.Item(363, 12).Value = 6000
$ws.Cells.Item(363, 13).Value = 5750
$ws.Cells.Item(363, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(363, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(363, 16).Value = 288
$ws.Cells.Item(363, 17).Value = 20
$ws.Cells.Item(363, 18).Value = "Hortaliza"
